$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the SMILES code for "AsymPolPOK" (row 14, column B)
$ws.Range("B14").Value = "`t[K+].[O-]P(=O)([O-].[K+])OC1CCC2(CC1)CC(NC(=O)C1=CC(C)(C)N([O])C1(C)C)CC1(CCC(OP(=O)([O-].[K+])O)CC1)N2[O]"

# Update the SMILES code for "AMUPolCbm" (row 13, column B)
$ws.Range("B13").Value = "`tN(CCOCOCOCOC)(C1CC2(CCOCC2)N([O])C2(CCOCC2)C1)C(OC1CC2(CCOCC2)N([O])C2(CCOCC2)C1)=O"

# Update the selected cell to B13
$ws.Range("B13").Select()

# Reposition the workbook window (xWindow/yWindow in the saved file)
$excel.ActiveWindow.Left = -42540
$excel.ActiveWindow.Top = 1580
